$d = $word.ActiveDocument

# Helper: do a single, scoped Find & Replace over the whole document content.
# Returns $true/$false depending on whether the search text was found.
function Replace-UniqueText($findText, $replaceText) {
    $rng = $d.Content
    return $rng.Find.Execute(
        $findText, $false, $false, $false, $false, $false,
        $true, 1, $false, $replaceText, 2)
}

# ---------------------------------------------------------------------------
# 1) Paragraph "đang diễn ra vô cùng ác liệt..." (answer C of Câu 1) gets the
#    "C. " prefix added.
# ---------------------------------------------------------------------------
$found1 = Replace-UniqueText `
    "đang diễn ra vô cùng ác liệt" `
    "C. đang diễn ra vô cùng ác liệt"
if (-not $found1) {
    Write-Host "WARNING: could not find text for edit 1"
}

# ---------------------------------------------------------------------------
# 2) Paragraph "bùng nổ và ngày càng lan rộng." (answer D of Câu 1) gets the
#    "D. " prefix added.
# ---------------------------------------------------------------------------
$found2 = Replace-UniqueText `
    "bùng nổ và ngày càng lan rộng" `
    "D. bùng nổ và ngày càng lan rộng"
if (-not $found2) {
    Write-Host "WARNING: could not find text for edit 2"
}

# ---------------------------------------------------------------------------
# 3) Lower-case "a" that starts the answer "a. Anh, Pháp, Mỹ." (answer A of
#    Câu 2) becomes upper-case "A". Only the leading character is touched so
#    that the distinct formatting kept on the neighbouring "Mỹ" run (it is
#    tagged with lang="vi-VN") is left completely intact.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$foundAnchor = $anchor.Find.Execute(
    "Anh, Pháp, Mỹ", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if ($foundAnchor) {
    [void]$anchor.Expand(4) # wdParagraph
    $firstChar = $d.Range($anchor.Start, $anchor.Start + 1)
    if ($firstChar.Text -eq "a") {
        $firstChar.Text = "A"
    } else {
        Write-Host "WARNING: unexpected first character for edit 3:" $firstChar.Text
    }
} else {
    Write-Host "WARNING: could not find anchor paragraph for edit 3"
}

# ---------------------------------------------------------------------------
# 4) Insert a missing space between "C." and "Anh, Pháp, Liên Xô" (answer C
#    of Câu 2).
# ---------------------------------------------------------------------------
$found4 = Replace-UniqueText `
    "C.Anh, Pháp, Liên Xô" `
    "C. Anh, Pháp, Liên Xô"
if (-not $found4) {
    Write-Host "WARNING: could not find text for edit 4"
}

# ---------------------------------------------------------------------------
# 5) Insert a missing space between "D." and "Liên Xô, Mỹ, Anh" (answer D of
#    Câu 2).
# ---------------------------------------------------------------------------
$found5 = Replace-UniqueText `
    "D.Liên Xô, " `
    "D. Liên Xô, "
if (-not $found5) {
    Write-Host "WARNING: could not find text for edit 5"
}
